$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price / 1h-volume snapshot (and restore the
# FraxShare/FTXToken row order at 46-47). Price cells that look like plain
# decimals get NumberFormat "@" first so Excel keeps them as text (matching
# the source data's leading/trailing zeros) instead of coercing to numbers.
$ws.Range("D2").Value = '38.636.97'
$ws.Range("E2").Value = '  +2.32%  '
$ws.Range("D3").Value = '2.095.00'
$ws.Range("E3").Value = '  +3.04%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.72'
$ws.Range("E5").Value = '  +0.31%  '
$ws.Range("E6").Value = '  +1.46%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '61.40'
$ws.Range("E7").Value = '  +2.04%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  +1.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0844'
$ws.Range("E10").Value = '  +2.58%  '
$ws.Range("E11").Value = '  +0.26%  '
$ws.Range("D12").Value = '2.401.64'
$ws.Range("E12").Value = '  +2.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.82'
$ws.Range("E13").Value = '  +1.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.36'
$ws.Range("E14").Value = '  +6.23%  '
$ws.Range("E15").Value = '  +1.46%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.46'
$ws.Range("E16").Value = '  +5.15%  '
$ws.Range("D17").Value = '2.093.91'
$ws.Range("E17").Value = '  +2.84%  '
$ws.Range("D18").Value = '38.575.10'
$ws.Range("E18").Value = '  +2.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.02'
$ws.Range("E19").Value = '  +2.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.06'
$ws.Range("E20").Value = '  +3.10%  '
$ws.Range("D21").Value = '0.0₃0836'
$ws.Range("E21").Value = '  +1.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '226.23'
$ws.Range("E22").Value = '  +1.09%  '
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.40'
$ws.Range("E24").Value = '  -0.34%  '
$ws.Range("E25").Value = '  +1.86%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '170.56'
$ws.Range("E26").Value = '  +1.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.45'
$ws.Range("E27").Value = '  +1.02%  '
$ws.Range("E28").Value = '  +2.90%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.12'
$ws.Range("E29").Value = '  +1.84%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.36'
$ws.Range("E30").Value = '  +7.42%  '
$ws.Range("E31").Value = '  +0.09%  '
$ws.Range("E32").Value = '  +3.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.81'
$ws.Range("E33").Value = '  +6.70%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.48'
$ws.Range("E34").Value = '  +2.40%  '
$ws.Range("E35").Value = '  +0.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.56'
$ws.Range("E36").Value = '  +2.00%  '
$ws.Range("E37").Value = '  +3.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.59'
$ws.Range("E38").Value = '  +4.88%  '
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.56'
$ws.Range("E40").Value = '  +2.32%  '
$ws.Range("D41").Value = '1.546.28'
$ws.Range("E41").Value = '  +0.64%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.83'
$ws.Range("E42").Value = '  +4.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0219'
$ws.Range("E43").Value = '  +1.72%  '
$ws.Range("E44").Value = '  +1.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0915'
$ws.Range("E45").Value = '  +0.43%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.70'
$ws.Range("E46").Value = '  +8.62%  '
$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.15'
$ws.Range("E47").Value = '  +3.24%  '
$ws.Range("E48").Value = '  +1.09%  '
$ws.Range("E49").Value = '  +2.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.00'
$ws.Range("E50").Value = '  +1.07%  '
$ws.Range("D51").Value = '2.292.76'
$ws.Range("E51").Value = '  +3.13%  '
